$wb = $excel.ActiveWorkbook

# --- Rename the "Estimated" sheet to "Initial Estimate" ---
$wsEstimated = $wb.Worksheets.Item("Estimated")
$wsEstimated.Name = "Initial Estimate"

$wsActual = $wb.Worksheets.Item("Actual")
$wsAnalysis = $wb.Worksheets.Item("Analysis")

# --- Analysis sheet: remove the lone "Grand Total" entry (now tracked on Actual) ---
$wsAnalysis.Range("A2").ClearContents() | Out-Null

# --- Actual sheet: rework the task list ---
# Clear the old date value/format out of A2, then turn it into the "Batch 1" label.
$wsActual.Range("A2").ClearFormats() | Out-Null
$wsActual.Range("A2").Value = "Batch 1"

# A3 keeps the old date-style formatting (now blank) while B3:B12 hold the batch items.
$wsActual.Range("A3").NumberFormat = "m/d/yy"

$wsActual.Range("B3").Value = "ABS Filament"
$wsActual.Range("B4").Value = "TPU Filament"
$wsActual.Range("B5").Value = "Cotton Wicks"
$wsActual.Range("B6").Value = "Lavender Essential Oil"
$wsActual.Range("B7").Value = "Differently Scented Essential Oil"
$wsActual.Range("B8").Value = "Silicone"
$wsActual.Range("B9").Value = "Mold Release Spray"
$wsActual.Range("B10").Value = "Valves"
$wsActual.Range("B11").Value = "Pumps"
$wsActual.Range("B12").Value = "Wires"

# Wrap the longer notes so they read well in a narrower column.
$wsActual.Range("B6").WrapText = $true
$wsActual.Range("B7").WrapText = $true
$wsActual.Range("B9").WrapText = $true

$wsActual.Rows.Item(6).RowHeight = 30
$wsActual.Rows.Item(7).RowHeight = 45
$wsActual.Rows.Item(9).RowHeight = 30.75

# Column B needs to be wide enough for the new item names.
$wsActual.Columns.Item(2).ColumnWidth = 11.6

# --- Selection / active-sheet bookkeeping, applied last so "Actual" ends up active ---
$wsAnalysis.Range("A2").Select() | Out-Null
$wsActual.Activate() | Out-Null
$wsActual.Range("B13").Select() | Out-Null
